$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet: Correspond Handoff Datetime (D) and Correspond Handback DateTime (G)
$wsZhCn.Range("D2:D3").Value = "2016-02-18 02:24:09"
$wsZhCn.Range("G2:G3").Value = "2016-02-18 02:24:54"

# de-de sheet: Correspond Handoff Datetime (D) and Correspond Handback DateTime (G)
$wsDeDe.Range("D2:D3").Value = "2016-02-18 02:24:20"
$wsDeDe.Range("G2:G3").Value = "2016-02-18 02:25:14"
